$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("description (รายละเอียด)")

$ws2.Range("G1:H1").Merge()
$ws2.Range("G1").HorizontalAlignment = -4108
$ws2.Columns("G:H").ColumnWidth = 25.7265625
